$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold plain decimal text (e.g. "0.9990") that Excel would
# otherwise auto-convert to a Number, dropping the trailing zero(s).
# Force them to Text format first so the literal string is preserved,
# matching the workbook author's original inlineStr text cells.
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "30.485.19"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.911.58"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "244.91"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +3.09%  "
$ws.Range("D8").Value = "0.2893"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "0.06701"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "110.55"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").Value = "19.07"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("D12").Value = "1.911.89"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "0.07549"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "5.265"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "0.6686"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "273.85"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "30.473.87"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "0.000007532"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "2.167.19"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "5.482"
$ws.Range("E22").Value = "  +5.45%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "6.451"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "9.440"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "163.45"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "20.17"
$ws.Range("E27").Value = "  -4.98%  "
$ws.Range("D28").Value = "2.115"
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("D29").Value = "0.1048"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "4.134"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "4.055"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").Value = "0.04989"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").Value = "0.7294"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.725"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02029"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "110.76"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").Value = "2.021"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "0.4417"
$ws.Range("E42").Value = "  +5.61%  "
$ws.Range("D43").Value = "0.8664"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "5.857"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "0.9991"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "67.90"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "7.411"
$ws.Range("E47").Value = "  +4.49%  "
$ws.Range("D48").Value = "9.189"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "0.1242"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("D50").Value = "47.57"
$ws.Range("E50").Value = "  -9.54%  "
$ws.Range("D51").Value = "1.469"
$ws.Range("E51").Value = "  +7.51%  "
